$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: B=Coin name, C=Link, D=Price (text), E=Volume/1h change (text)
# D-column price values that parse as plain numbers must be force-written as
# text (matching the source data, which stores every Price cell as a string)
# without leaving a residual NumberFormat/quote-prefix style on the cell.

$ws.Range('D2').Value = '26.504.49'
$ws.Range('E2').Value = '  -1.65%  '
$ws.Range('D3').Value = '1.790.66'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '303.20'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4251'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.73%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3616'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07167'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8504'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.59'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.46%  '
$ws.Range('D12').Value = '1.802.32'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.492'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.266'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06910'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.007'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.40'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008768'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.96'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').Value = '26.501.72'
$ws.Range('E21').Value = '  -1.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.134'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('D24').Value = '1.997.91'
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.43'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.812'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.14'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.123'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.867'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +14.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08885'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7391'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.138'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.353'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.746'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.111'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05161'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01896'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4965'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.71%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1619'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.613'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.379'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.192'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.32'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '105.31'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.002'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.636'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.08%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4513'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06209'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.765'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.48%  '
